$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "skills" column before the old "start.month" column (G) ---
$ws.Columns("G:G").Insert()

# Header for the new column
$ws.Range("G1").Value2 = "skills"

# Match the visual formatting of the new column's data cells to the
# "description" column (wrapped text), like the real workbook shows.
$ws.Range("D2").Copy()
$ws.Range("G2:G5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the "description" column with the real project descriptions ---
$ws.Range("D2").Value2 = "Create and automate an ETL processing pipeline for COVID-19 data using Python and AWS. Infrastructure was created and managed by CloudFormation. Data was fetched from a hosted GitHub repositoy on a schedule managed by EventBridge and transformed using Python code running within a Lambda function. Transformed data was stored in a DynamoDB database and streamed to an S3 bucket. QuickSight was connected to the S3 bucket so COVID-19 data could be displayed on an easy to read dashboard. Code was tracked by Git version control and stored in a GitHub repository."
$ws.Range("D3").Value2 = "Create a Netflix style recommendation engine using SageMaker. Data was loaded from IMDB and uploaded to S3 using Jupyter Notebooks. Jupyter was also used to conduct the machine learning training. Feature engineering was done with Athena and SQL. Machine learning results were exported and used to build a PHP website. The website was fully hosted on AWS. DNS registration and routing was done with Route53. Dynamic content was served by a containerized PHP Lambda using Bref fronted by API Gateway. Static content was served by S3. All content was distributed using CloudFront."
$ws.Range("D4").Value2 = "Create an image processing website that utilizes multiple cloud providers. Used AWS, Azure, and GCP to accomplish the goal. For AWS, Route 53 was used for DNS registration and routing, CloudFront for content distribution, API Gateway for HTTP access to Lambda functions, C# and Node.js Lambda functions for business logic, and S3 for website hosting. For Azure, Table Storage was used to store image analysis data in a NoSQL database. For GCP, Cloud Vision was used to conduct image analysis on images uploaded to S3."
$ws.Range("D5").Value2 = "Create a verification website in Azure that is continously integrated and globally performant. The website was written in Ruby on Rails and deployed to App Service using Azure Pipelines. The app servers were setup with auto scaling and secured within a Virtual Network behind Front Door. Data was written to Cosmos DB and binary data was written to Blob storage. Azure Pipelines was executed using a self-hosted agent running on a Virtual Machine with a local SSD disk and public IP address protected by a Network Security Group. All infrastructure was deployed using Resource Manager."

# --- Fill in the new "skills" data column ---
$ws.Range("G2").Value2 = "CloudFormation, DynamoDB,Event Bridge,Git,Lambda,Python,QuickSight,Serverless,Simple Notifcation Service (SNS),Simple Storage Service (S3)"
$ws.Range("G3").Value2 = "API Gateway,Athena,Bref,CloudFront,Elastic Container Registry (ECR),Git,Jupyter Notebooks,Lambda,Machine Learning,PHP,Python,Route 53,SageMaker,Serverless,Simple Storage Service (S3),SQL"
$ws.Range("G4").Value2 = "API Gateway,C#,Cloud Vision,CloudFront,Git,Lambda,Node.js,Route 53,Serverless,Simple Storage Service (S3),Table Storage"
$ws.Range("G5").Value2 = "App Service,Blob,CosmosDB,Disk Storage,Front Door,Git,MongoDB,Network Security Group (NSG),Pipelines,Public IP,Rails,Resource Manager,Ruby,Virtual Machine,Virtual Network"

# Give the new column its own (narrower) width, distinct from the bestFit columns.
$ws.Columns("G:G").ColumnWidth = 9.8

# --- Leftover empty-cell artifacts from editing below the table (as in the source file) ---
$ws.Range("C11").Value2 = "x"
$ws.Range("C11").ClearContents()
$ws.Range("C11").Style = "Normal"
$ws.Rows(11).RowHeight = 15

$ws.Range("E18").Value2 = "x"
$ws.Range("E18").ClearContents()
$ws.Range("E18").Style = "Normal"
$ws.Rows(18).RowHeight = 15

# --- View state: zoom + selection ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("G3").Select()
